$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B28 should become a true numeric value (it currently is stored as text "4")
$ws.Range("B28").Value = 4

# Add new row 29 data
$ws.Range("A29").Value = "Ruilin"

# B29 holds "3" but must remain text (not auto-converted to a number)
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "3"
$ws.Range("B29").ClearFormats()

$ws.Range("C29").Value = "无"
$ws.Range("D29").Value = "FBK"
$ws.Range("E29").Value = "WRI"
$ws.Range("F29").Value = "9b81a0cf-ae6f-4476-b619-1b75e1becf94"
$ws.Range("G29").Value = "B1ae1lZRb_annotated.xlsx"
$ws.Range("H29").Value = "We will elaborate on this aspect in the final version of the paper."
